# Update vm_pu.xlsx results for the 380 kV case (Case_0_158)
# Sets bus voltage magnitude (p.u.) values for rows 2-25, columns B-F and I-N
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.047039682708714
$ws.Cells.Item(2, 4).Value = 1.049684267066092
$ws.Cells.Item(2, 5).Value = 1.053884458289688
$ws.Cells.Item(2, 6).Value = 1.063063112619902
$ws.Cells.Item(2, 9).Value = 1.04149200229736
$ws.Cells.Item(2, 10).Value = 1.052090871787959
$ws.Cells.Item(2, 11).Value = 1.052440640327322
$ws.Cells.Item(2, 12).Value = 1.056629195115558
$ws.Cells.Item(2, 13).Value = 1.065782766879213
$ws.Cells.Item(2, 14).Value = 1.021131664836827
# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.048330387718712
$ws.Cells.Item(3, 4).Value = 1.050890545911572
$ws.Cells.Item(3, 5).Value = 1.055035962485241
$ws.Cells.Item(3, 6).Value = 1.064311590246407
$ws.Cells.Item(3, 9).Value = 1.041792396837425
$ws.Cells.Item(3, 10).Value = 1.053028209899666
$ws.Cells.Item(3, 11).Value = 1.053458197933901
$ws.Cells.Item(3, 12).Value = 1.057592966904132
$ws.Cells.Item(3, 13).Value = 1.066845093708552
$ws.Cells.Item(3, 14).Value = 1.021449335311988
# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.049164757177336
$ws.Cells.Item(4, 4).Value = 1.051670573656668
$ws.Cells.Item(4, 5).Value = 1.055780660943154
$ws.Cells.Item(4, 6).Value = 1.065119045730061
$ws.Cells.Item(4, 9).Value = 1.041984569936853
$ws.Cells.Item(4, 10).Value = 1.053633443332907
$ws.Cells.Item(4, 11).Value = 1.054115544499963
$ws.Cells.Item(4, 12).Value = 1.05821561051617
$ws.Cells.Item(4, 13).Value = 1.067531535133883
$ws.Cells.Item(4, 14).Value = 1.02165429115508
# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.049515337253012
$ws.Cells.Item(5, 4).Value = 1.051998376500986
$ws.Cells.Item(5, 5).Value = 1.056093638270767
$ws.Cells.Item(5, 6).Value = 1.065458408457293
$ws.Cells.Item(5, 9).Value = 1.042064833022041
$ws.Cells.Item(5, 10).Value = 1.053887577522354
$ws.Cells.Item(5, 11).Value = 1.054391636429154
$ws.Cells.Item(5, 12).Value = 1.058477137090774
$ws.Cells.Item(5, 13).Value = 1.067819888722743
$ws.Cells.Item(5, 14).Value = 1.02174031204181
# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.049574190214718
$ws.Cells.Item(6, 4).Value = 1.052053409029891
$ws.Cells.Item(6, 5).Value = 1.056146183100581
$ws.Cells.Item(6, 6).Value = 1.065515383656783
$ws.Cells.Item(6, 9).Value = 1.042078278695095
$ws.Cells.Item(6, 10).Value = 1.053930229892851
$ws.Cells.Item(6, 11).Value = 1.05443797849986
$ws.Cells.Item(6, 12).Value = 1.058521034962615
$ws.Cells.Item(6, 13).Value = 1.067868291349106
$ws.Cells.Item(6, 14).Value = 1.021754746987111
# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.049169442383289
$ws.Cells.Item(7, 4).Value = 1.0516749542459
$ws.Cells.Item(7, 5).Value = 1.055784843327686
$ws.Cells.Item(7, 6).Value = 1.065123580668059
$ws.Cells.Item(7, 9).Value = 1.041985644483959
$ws.Cells.Item(7, 10).Value = 1.053636840285873
$ws.Cells.Item(7, 11).Value = 1.054119234659146
$ws.Cells.Item(7, 12).Value = 1.058219105959216
$ws.Cells.Item(7, 13).Value = 1.067535389015406
$ws.Cells.Item(7, 14).Value = 1.021655441129687
# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.04747605082097
$ws.Cells.Item(8, 4).Value = 1.050092042892357
$ws.Cells.Item(8, 5).Value = 1.054273699084897
$ws.Cells.Item(8, 6).Value = 1.063485124386476
$ws.Cells.Item(8, 9).Value = 1.041593978477869
$ws.Cells.Item(8, 10).Value = 1.052407917381639
$ws.Cells.Item(8, 11).Value = 1.052784753944409
$ws.Cells.Item(8, 12).Value = 1.05695511073002
$ws.Cells.Item(8, 13).Value = 1.066141984322726
$ws.Cells.Item(8, 14).Value = 1.021239147320228
# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.044485775918684
$ws.Cells.Item(9, 4).Value = 1.047298688905325
$ws.Cells.Item(9, 5).Value = 1.051607690071559
$ws.Cells.Item(9, 6).Value = 1.060594830664033
$ws.Cells.Item(9, 9).Value = 1.040886912487681
$ws.Cells.Item(9, 10).Value = 1.05023244348836
$ws.Cells.Item(9, 11).Value = 1.050424844137705
$ws.Cells.Item(9, 12).Value = 1.054720180268962
$ws.Cells.Item(9, 13).Value = 1.063679210472267
$ws.Cells.Item(9, 14).Value = 1.0205009731277
# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.042487805689019
$ws.Cells.Item(10, 4).Value = 1.045433557923191
$ws.Cells.Item(10, 5).Value = 1.049828066141052
$ws.Cells.Item(10, 6).Value = 1.058665704560969
$ws.Cells.Item(10, 9).Value = 1.040404122951353
$ws.Cells.Item(10, 10).Value = 1.048775295073274
$ws.Cells.Item(10, 11).Value = 1.048845786670775
$ws.Cells.Item(10, 12).Value = 1.053224975442699
$ws.Cells.Item(10, 13).Value = 1.062032228980429
$ws.Cells.Item(10, 14).Value = 1.020005713861591
# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.041621557283808
$ws.Cells.Item(11, 4).Value = 1.044625211562318
$ws.Cells.Item(11, 5).Value = 1.049056893926014
$ws.Cells.Item(11, 6).Value = 1.057829797556887
$ws.Cells.Item(11, 9).Value = 1.04019235131929
$ws.Cells.Item(11, 10).Value = 1.048142681269283
$ws.Cells.Item(11, 11).Value = 1.048160631798496
$ws.Cells.Item(11, 12).Value = 1.052576258859208
$ws.Cells.Item(11, 13).Value = 1.061317817194657
$ws.Cells.Item(11, 14).Value = 1.019790505788593
# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.041299622357306
$ws.Cells.Item(12, 4).Value = 1.044324842376392
$ws.Cells.Item(12, 5).Value = 1.048770355600059
$ws.Cells.Item(12, 6).Value = 1.05751921362156
$ws.Cells.Item(12, 9).Value = 1.040113280081017
$ws.Cells.Item(12, 10).Value = 1.047907448247897
$ws.Cells.Item(12, 11).Value = 1.047905919547128
$ws.Cells.Item(12, 12).Value = 1.052335101114047
$ws.Cells.Item(12, 13).Value = 1.061052260877911
$ws.Cells.Item(12, 14).Value = 1.019710453226932
# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.041368686331201
$ws.Cells.Item(13, 4).Value = 1.044389277769718
$ws.Cells.Item(13, 5).Value = 1.048831823193091
$ws.Cells.Item(13, 6).Value = 1.057585839094869
$ws.Cells.Item(13, 9).Value = 1.040130259693105
$ws.Cells.Item(13, 10).Value = 1.047957917983019
$ws.Cells.Item(13, 11).Value = 1.04796056597753
$ws.Cells.Item(13, 12).Value = 1.052386839161348
$ws.Cells.Item(13, 13).Value = 1.061109232344455
$ws.Cells.Item(13, 14).Value = 1.019727629977969
# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.041594949570021
$ws.Cells.Item(14, 4).Value = 1.044600385283479
$ws.Cells.Item(14, 5).Value = 1.049033210444794
$ws.Cells.Item(14, 6).Value = 1.057804126464585
$ws.Cells.Item(14, 9).Value = 1.040185823632717
$ws.Cells.Item(14, 10).Value = 1.048123241993436
$ws.Cells.Item(14, 11).Value = 1.048139581612315
$ws.Cells.Item(14, 12).Value = 1.052556328689395
$ws.Cells.Item(14, 13).Value = 1.061295870160126
$ws.Cells.Item(14, 14).Value = 1.019783890963908
# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.041734334937733
$ws.Cells.Item(15, 4).Value = 1.044730440469355
$ws.Cells.Item(15, 5).Value = 1.04915727967808
$ws.Cells.Item(15, 6).Value = 1.057938608378677
$ws.Cells.Item(15, 9).Value = 1.04022000407117
$ws.Cells.Item(15, 10).Value = 1.048225070082905
$ws.Cells.Item(15, 11).Value = 1.048249850443532
$ws.Cells.Item(15, 12).Value = 1.05266073078998
$ws.Cells.Item(15, 13).Value = 1.061410838350726
$ws.Cells.Item(15, 14).Value = 1.019818539988464
# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.042545271901857
$ws.Cells.Item(16, 4).Value = 1.045487189495692
$ws.Cells.Item(16, 5).Value = 1.049879233756517
$ws.Cells.Item(16, 6).Value = 1.058721168397167
$ws.Cells.Item(16, 9).Value = 1.040418120125023
$ws.Cells.Item(16, 10).Value = 1.048817244373326
$ws.Cells.Item(16, 11).Value = 1.048891228126136
$ws.Cells.Item(16, 12).Value = 1.053268001358678
$ws.Cells.Item(16, 13).Value = 1.062079615380737
$ws.Cells.Item(16, 14).Value = 1.020019980476343
# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.043053649915662
$ws.Cells.Item(17, 4).Value = 1.045961679754886
$ws.Cells.Item(17, 5).Value = 1.050331938200728
$ws.Cells.Item(17, 6).Value = 1.059211889290323
$ws.Cells.Item(17, 9).Value = 1.040541663950662
$ws.Cells.Item(17, 10).Value = 1.049188253511575
$ws.Cells.Item(17, 11).Value = 1.049293167103103
$ws.Cells.Item(17, 12).Value = 1.053648580695911
$ws.Cells.Item(17, 13).Value = 1.062498782661184
$ws.Cells.Item(17, 14).Value = 1.020146135331019
# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.04335007100842
$ws.Cells.Item(18, 4).Value = 1.046238371462692
$ws.Cells.Item(18, 5).Value = 1.050595936967401
$ws.Cells.Item(18, 6).Value = 1.059498062492004
$ws.Cells.Item(18, 9).Value = 1.040613462552569
$ws.Cells.Item(18, 10).Value = 1.049404496997547
$ws.Cells.Item(18, 11).Value = 1.049527475220676
$ws.Cells.Item(18, 12).Value = 1.053870442579499
$ws.Cells.Item(18, 13).Value = 1.062743154675526
$ws.Cells.Item(18, 14).Value = 1.020219646292273
# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.04345112488735
$ws.Cells.Item(19, 4).Value = 1.046332704338047
$ws.Cells.Item(19, 5).Value = 1.050685944302737
$ws.Cells.Item(19, 6).Value = 1.059595630760496
$ws.Cells.Item(19, 9).Value = 1.040637899541682
$ws.Cells.Item(19, 10).Value = 1.049478203407343
$ws.Cells.Item(19, 11).Value = 1.049607345241017
$ws.Cells.Item(19, 12).Value = 1.053946070865693
$ws.Cells.Item(19, 13).Value = 1.062826458754382
$ws.Cells.Item(19, 14).Value = 1.020244699269672
# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.042999116887032
$ws.Cells.Item(20, 4).Value = 1.045910778751127
$ws.Cells.Item(20, 5).Value = 1.050283373121297
$ws.Cells.Item(20, 6).Value = 1.059159245399496
$ws.Cells.Item(20, 9).Value = 1.04052843601924
$ws.Cells.Item(20, 10).Value = 1.049148464299781
$ws.Cells.Item(20, 11).Value = 1.0492500569538
$ws.Cells.Item(20, 12).Value = 1.053607760940869
$ws.Cells.Item(20, 13).Value = 1.062453822534701
$ws.Cells.Item(20, 14).Value = 1.020132607669782
# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.041528325460179
$ws.Cells.Item(21, 4).Value = 1.044538222538433
$ws.Cells.Item(21, 5).Value = 1.048973909444558
$ws.Cells.Item(21, 6).Value = 1.057739848809666
$ws.Cells.Item(21, 9).Value = 1.040169472759381
$ws.Cells.Item(21, 10).Value = 1.048074565174823
$ws.Cells.Item(21, 11).Value = 1.048086871935212
$ws.Cells.Item(21, 12).Value = 1.052506423667441
$ws.Cells.Item(21, 13).Value = 1.061240915298601
$ws.Cells.Item(21, 14).Value = 1.019767326679304
# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.040602584960827
$ws.Cells.Item(22, 4).Value = 1.043674584462768
$ws.Cells.Item(22, 5).Value = 1.04815007181034
$ws.Cells.Item(22, 6).Value = 1.056846890518824
$ws.Cells.Item(22, 9).Value = 1.039941406641278
$ws.Cells.Item(22, 10).Value = 1.047397902244048
$ws.Cells.Item(22, 11).Value = 1.04735428461768
$ws.Cells.Item(22, 12).Value = 1.051812836600212
$ws.Cells.Item(22, 13).Value = 1.060477200206394
$ws.Cells.Item(22, 14).Value = 1.019536995951866
# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.041093433329058
$ws.Cells.Item(23, 4).Value = 1.044132478788685
$ws.Cells.Item(23, 5).Value = 1.048586854576296
$ws.Cells.Item(23, 6).Value = 1.05732031570016
$ws.Cells.Item(23, 9).Value = 1.040062533979075
$ws.Cells.Item(23, 10).Value = 1.047756753408106
$ws.Cells.Item(23, 11).Value = 1.04774276219879
$ws.Cells.Item(23, 12).Value = 1.052180628512682
$ws.Cells.Item(23, 13).Value = 1.060882166447926
$ws.Cells.Item(23, 14).Value = 1.019659161873004
# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.043023758346367
$ws.Cells.Item(24, 4).Value = 1.045933778941246
$ws.Cells.Item(24, 5).Value = 1.050305317763429
$ws.Cells.Item(24, 6).Value = 1.059183033080654
$ws.Cells.Item(24, 9).Value = 1.040534413962879
$ws.Cells.Item(24, 10).Value = 1.049166443825753
$ws.Cells.Item(24, 11).Value = 1.049269536995006
$ws.Cells.Item(24, 12).Value = 1.053626206012944
$ws.Cells.Item(24, 13).Value = 1.062474138454879
$ws.Cells.Item(24, 14).Value = 1.02013872046324
# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.045259602719513
$ws.Cells.Item(25, 4).Value = 1.048021336034137
$ws.Cells.Item(25, 5).Value = 1.052297309158129
$ws.Cells.Item(25, 6).Value = 1.061342427830212
$ws.Cells.Item(25, 9).Value = 1.041071713814157
$ws.Cells.Item(25, 10).Value = 1.050796049613497
$ws.Cells.Item(25, 11).Value = 1.051035946048011
$ws.Cells.Item(25, 12).Value = 1.055298878955197
$ws.Cells.Item(25, 13).Value = 1.064316790391569
$ws.Cells.Item(25, 14).Value = 1.020692359694427

Write-Host "Updated 264 cells in vm_pu sheet"
